$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.1752021563342318
$ws.Range("C2").Value = 0.5849056603773585
$ws.Range("J2").Value = 0.005390835579514825
$ws.Range("P2").Value = 0.137466307277628
$ws.Range("S2").Value = 0.09703504043126684
$ws.Range("C3").Value = 0.009174311926605505
$ws.Range("J3").Value = 0.02752293577981652
$ws.Range("P3").Value = 0.7431192660550459
$ws.Range("S3").Value = 0.2201834862385321
$ws.Range("J4").Value = 0.05128205128205128
$ws.Range("P4").Value = 0.7692307692307693
$ws.Range("S4").Value = 0.1794871794871795
$ws.Range("B6").Value = 0.08403361344537816
$ws.Range("D6").Value = 0.01680672268907563
$ws.Range("F6").Value = 0.06722689075630252
$ws.Range("J6").Value = 0.2352941176470588
$ws.Range("O6").Value = 0.01260504201680672
$ws.Range("Q6").Value = 0.1554621848739496
$ws.Range("R6").Value = 0.0546218487394958
$ws.Range("S6").Value = 0.3739495798319328
$ws.Range("B7").Value = 0.1025641025641026
$ws.Range("D7").Value = 0.01282051282051282
$ws.Range("F7").Value = 0.05982905982905983
$ws.Range("J7").Value = 0.1324786324786325
$ws.Range("O7").Value = 0.0170940170940171
$ws.Range("Q7").Value = 0.1965811965811966
$ws.Range("R7").Value = 0.07692307692307693
$ws.Range("S7").Value = 0.4017094017094017
$ws.Range("B8").Value = 0.1241258741258741
$ws.Range("D8").Value = 0.005244755244755245
$ws.Range("F8").Value = 0.06818181818181818
$ws.Range("J8").Value = 0.1101398601398601
$ws.Range("O8").Value = 0.01573426573426574
$ws.Range("Q8").Value = 0.1660839160839161
$ws.Range("R8").Value = 0.09265734265734266
$ws.Range("S8").Value = 0.4178321678321678
$ws.Range("B9").Value = 0.1417322834645669
$ws.Range("D9").Value = 0.02362204724409449
$ws.Range("F9").Value = 0.07874015748031496
$ws.Range("J9").Value = 0.07874015748031496
$ws.Range("O9").Value = 0.03149606299212598
$ws.Range("Q9").Value = 0.1102362204724409
$ws.Range("R9").Value = 0.1023622047244094
$ws.Range("S9").Value = 0.4330708661417323
$ws.Range("B10").Value = 0.1164574616457462
$ws.Range("D10").Value = 0.02092050209205021
$ws.Range("F10").Value = 0.06136680613668061
$ws.Range("J10").Value = 0.1408647140864714
$ws.Range("O10").Value = 0.0104602510460251
$ws.Range("Q10").Value = 0.2231520223152022
$ws.Range("R10").Value = 0.06276150627615062
$ws.Range("S10").Value = 0.3640167364016736
$ws.Range("G11").Value = 0.1598915989159892
$ws.Range("J11").Value = 0.1084010840108401
$ws.Range("K11").Value = 0.2330623306233062
$ws.Range("L11").Value = 0.4905149051490515
$ws.Range("S11").Value = 0.008130081300813009
$ws.Range("G12").Value = 0.7864583333333334
$ws.Range("J12").Value = 0.1302083333333333
$ws.Range("K12").Value = 0.005208333333333333
$ws.Range("L12").Value = 0.04166666666666666
$ws.Range("S12").Value = 0.03645833333333334
$ws.Range("G13").Value = 0.4363636363636363
$ws.Range("J13").Value = 0.4545454545454545
$ws.Range("S13").Value = 0.1090909090909091
$ws.Range("F15").Value = 0.008888888888888889
$ws.Range("H15").Value = 0.2088888888888889
$ws.Range("I15").Value = 0.01777777777777778
$ws.Range("J15").Value = 0.3955555555555555
$ws.Range("K15").Value = 0.07111111111111111
$ws.Range("M15").Value = 0.008888888888888889
$ws.Range("O15").Value = 0.03555555555555556
$ws.Range("S15").Value = 0.2533333333333334
$ws.Range("F16").Value = 0.02127659574468085
$ws.Range("H16").Value = 0.148936170212766
$ws.Range("I16").Value = 0.06382978723404255
$ws.Range("J16").Value = 0.4212765957446808
$ws.Range("K16").Value = 0.1106382978723404
$ws.Range("M16").Value = 0.03404255319148936
$ws.Range("N16").Value = 0.00425531914893617
$ws.Range("O16").Value = 0.03829787234042553
$ws.Range("S16").Value = 0.1574468085106383
$ws.Range("F17").Value = 0.01972386587771203
$ws.Range("H17").Value = 0.2149901380670611
$ws.Range("I17").Value = 0.04930966469428008
$ws.Range("J17").Value = 0.4201183431952663
$ws.Range("K17").Value = 0.09467455621301775
$ws.Range("M17").Value = 0.01577909270216963
$ws.Range("N17").Value = 0.003944773175542407
$ws.Range("O17").Value = 0.07889546351084813
$ws.Range("S17").Value = 0.1025641025641026
$ws.Range("F18").Value = 0.02185792349726776
$ws.Range("H18").Value = 0.1748633879781421
$ws.Range("I18").Value = 0.06557377049180328
$ws.Range("J18").Value = 0.4426229508196721
$ws.Range("K18").Value = 0.08743169398907104
$ws.Range("M18").Value = 0.02185792349726776
$ws.Range("O18").Value = 0.07103825136612021
$ws.Range("S18").Value = 0.1147540983606557
$ws.Range("F19").Value = 0.01394700139470014
$ws.Range("H19").Value = 0.2412831241283124
$ws.Range("I19").Value = 0.04951185495118549
$ws.Range("J19").Value = 0.3619246861924686
$ws.Range("K19").Value = 0.1220362622036262
$ws.Range("M19").Value = 0.02301255230125523
$ws.Range("N19").Value = 0.002789400278940028
$ws.Range("O19").Value = 0.06066945606694561
$ws.Range("S19").Value = 0.1248256624825663
